$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CSR user name in A3 from "demosalesCSR" to "DemoCSR"
$ws.Range("A3").Value = "DemoCSR"

# Reflect the new active selection left behind on the sheet after the edit
$ws.Range("A3").Select()
